# "Started statistics per country": collapse the four side-by-side
# (Country, Code) tables in all_codes (A:B, C:D, E:F, G:H) into a single
# tall (country, code) table in columns A:B so it can be used as a lookup.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("all_codes")

# Drop the three extra side-by-side tables; A:B keeps the first one for now,
# it gets overwritten below along with the header.
$ws1.Range("C1:H14").Clear()

# New lower-case headers for the single lookup table.
$ws1.Range("A1").Value2 = "country"
$ws1.Range("B1").Value2 = "code"

# Flattened (country, code) pairs - these used to be read down column A then C
# then E then G (and B/D/F/H for the codes); now they all live in column A/B.
$ws1.Cells.Item(2, 1).Value2 = "Albania"
$ws1.Cells.Item(2, 2).Value2 = "al"
$ws1.Cells.Item(3, 1).Value2 = "Andorra"
$ws1.Cells.Item(3, 2).Value2 = "ad"
$ws1.Cells.Item(4, 1).Value2 = "Armenia"
$ws1.Cells.Item(4, 2).Value2 = "am"
$ws1.Cells.Item(5, 1).Value2 = "Australia"
$ws1.Cells.Item(5, 2).Value2 = "au"
$ws1.Cells.Item(6, 1).Value2 = "Austria"
$ws1.Cells.Item(6, 2).Value2 = "at"
$ws1.Cells.Item(7, 1).Value2 = "Azerbaijan"
$ws1.Cells.Item(7, 2).Value2 = "az"
$ws1.Cells.Item(8, 1).Value2 = "Belarus"
$ws1.Cells.Item(8, 2).Value2 = "by"
$ws1.Cells.Item(9, 1).Value2 = "Belgium"
$ws1.Cells.Item(9, 2).Value2 = "be"
$ws1.Cells.Item(10, 1).Value2 = "Bosnia and Herzegovina"
$ws1.Cells.Item(10, 2).Value2 = "ba"
$ws1.Cells.Item(11, 1).Value2 = "Bulgaria"
$ws1.Cells.Item(11, 2).Value2 = "bg"
$ws1.Cells.Item(12, 1).Value2 = "Croatia"
$ws1.Cells.Item(12, 2).Value2 = "hr"
$ws1.Cells.Item(13, 1).Value2 = "Cyprus"
$ws1.Cells.Item(13, 2).Value2 = "cy"
$ws1.Cells.Item(14, 1).Value2 = "Czechia"
$ws1.Cells.Item(14, 2).Value2 = "cz"
$ws1.Cells.Item(15, 1).Value2 = "Denmark"
$ws1.Cells.Item(15, 2).Value2 = "dk"
$ws1.Cells.Item(16, 1).Value2 = "Estonia"
$ws1.Cells.Item(16, 2).Value2 = "ee"
$ws1.Cells.Item(17, 1).Value2 = "Finland"
$ws1.Cells.Item(17, 2).Value2 = "fi"
$ws1.Cells.Item(18, 1).Value2 = "France"
$ws1.Cells.Item(18, 2).Value2 = "fr"
$ws1.Cells.Item(19, 1).Value2 = "Georgia"
$ws1.Cells.Item(19, 2).Value2 = "ge"
$ws1.Cells.Item(20, 1).Value2 = "Germany"
$ws1.Cells.Item(20, 2).Value2 = "de"
$ws1.Cells.Item(21, 1).Value2 = "Greece"
$ws1.Cells.Item(21, 2).Value2 = "gr"
$ws1.Cells.Item(22, 1).Value2 = "Hungary"
$ws1.Cells.Item(22, 2).Value2 = "hu"
$ws1.Cells.Item(23, 1).Value2 = "Iceland"
$ws1.Cells.Item(23, 2).Value2 = "is"
$ws1.Cells.Item(24, 1).Value2 = "Ireland"
$ws1.Cells.Item(24, 2).Value2 = "ie"
$ws1.Cells.Item(25, 1).Value2 = "Israel"
$ws1.Cells.Item(25, 2).Value2 = "il"
$ws1.Cells.Item(26, 1).Value2 = "Italy"
$ws1.Cells.Item(26, 2).Value2 = "it"
$ws1.Cells.Item(27, 1).Value2 = "Latvia"
$ws1.Cells.Item(27, 2).Value2 = "lv"
$ws1.Cells.Item(28, 1).Value2 = "Lithuania"
$ws1.Cells.Item(28, 2).Value2 = "lt"
$ws1.Cells.Item(29, 1).Value2 = "Luxembourg"
$ws1.Cells.Item(29, 2).Value2 = "lu"
$ws1.Cells.Item(30, 1).Value2 = "Malta"
$ws1.Cells.Item(30, 2).Value2 = "mt"
$ws1.Cells.Item(31, 1).Value2 = "Moldova"
$ws1.Cells.Item(31, 2).Value2 = "md"
$ws1.Cells.Item(32, 1).Value2 = "Monaco"
$ws1.Cells.Item(32, 2).Value2 = "mc"
$ws1.Cells.Item(33, 1).Value2 = "Montenegro"
$ws1.Cells.Item(33, 2).Value2 = "me"
$ws1.Cells.Item(34, 1).Value2 = "Morocco"
$ws1.Cells.Item(34, 2).Value2 = "ma"
$ws1.Cells.Item(35, 1).Value2 = "Netherlands"
$ws1.Cells.Item(35, 2).Value2 = "nl"
$ws1.Cells.Item(36, 1).Value2 = "North Macedonia"
$ws1.Cells.Item(36, 2).Value2 = "mk"
$ws1.Cells.Item(37, 1).Value2 = "Norway"
$ws1.Cells.Item(37, 2).Value2 = "no"
$ws1.Cells.Item(38, 1).Value2 = "Poland"
$ws1.Cells.Item(38, 2).Value2 = "pl"
$ws1.Cells.Item(39, 1).Value2 = "Portugal"
$ws1.Cells.Item(39, 2).Value2 = "pt"
$ws1.Cells.Item(40, 1).Value2 = "Romania"
$ws1.Cells.Item(40, 2).Value2 = "ro"
$ws1.Cells.Item(41, 1).Value2 = "Russia"
$ws1.Cells.Item(41, 2).Value2 = "ru"
$ws1.Cells.Item(42, 1).Value2 = "San Marino"
$ws1.Cells.Item(42, 2).Value2 = "sm"
$ws1.Cells.Item(43, 1).Value2 = "Serbia"
$ws1.Cells.Item(43, 2).Value2 = "rs"
$ws1.Cells.Item(44, 1).Value2 = "Serbia and Montenegro"
$ws1.Cells.Item(44, 2).Value2 = "cs"
$ws1.Cells.Item(45, 1).Value2 = "Slovakia"
$ws1.Cells.Item(45, 2).Value2 = "sk"
$ws1.Cells.Item(46, 1).Value2 = "Slovenia"
$ws1.Cells.Item(46, 2).Value2 = "si"
$ws1.Cells.Item(47, 1).Value2 = "Spain"
$ws1.Cells.Item(47, 2).Value2 = "es"
$ws1.Cells.Item(48, 1).Value2 = "Sweden"
$ws1.Cells.Item(48, 2).Value2 = "se"
$ws1.Cells.Item(49, 1).Value2 = "Switzerland"
$ws1.Cells.Item(49, 2).Value2 = "ch"
$ws1.Cells.Item(50, 1).Value2 = "Türkiye"
$ws1.Cells.Item(50, 2).Value2 = "tr"
$ws1.Cells.Item(51, 1).Value2 = "Ukraine"
$ws1.Cells.Item(51, 2).Value2 = "ua"
$ws1.Cells.Item(52, 1).Value2 = "United Kingdom"
$ws1.Cells.Item(52, 2).Value2 = "gb"
$ws1.Cells.Item(53, 1).Value2 = "Yugoslavia"
$ws1.Cells.Item(53, 2).Value2 = "yu"

# all_codes becomes the active sheet/tab, selection parked on the first code cell.
[void]$ws1.Activate()
[void]$ws1.Range("B2").Select()
